$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# ALC row 112: Making Ends Meet | Superior Spiritbond Potion
$ws_ALC.Range("H112").Value = 1649.5
$ws_ALC.Range("I112").Value = 1866
$ws_ALC.Range("J112").Value = 1000
$ws_ALC.Range("K112").Value = 5598
$ws_ALC.Range("L112").Value = 3000
$ws_ALC.Range("M112").Value = -4490
$ws_ALC.Range("N112").Value = -5216

$ws_ARM = $wb.Worksheets.Item("ARM")
# ARM row 45: Hollow Hallmarks | Mythril Ingot
$ws_ARM.Range("H45").Value = 15476.875
$ws_ARM.Range("I45").Value = 15476.875
$ws_ARM.Range("J45").Value = 0
$ws_ARM.Range("K45").Value = 15476.875
$ws_ARM.Range("L45").Value = 0
$ws_ARM.Range("M45").Value = -15099.875
$ws_ARM.Range("N45").ClearContents()

# ARM row 74: As the Bolt Flies | Titanium Nugget
$ws_ARM.Range("H74").Value = 7565.1113
$ws_ARM.Range("I74").Value = 4049.5
$ws_ARM.Range("J74").Value = 8569.571
$ws_ARM.Range("K74").Value = 4049.5
$ws_ARM.Range("L74").Value = 8569.571
$ws_ARM.Range("M74").Value = -3175.5
$ws_ARM.Range("N74").Value = -10317.571

# ARM row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws_ARM.Range("H77").Value = 7565.1113
$ws_ARM.Range("I77").Value = 4049.5
$ws_ARM.Range("J77").Value = 8569.571
$ws_ARM.Range("K77").Value = 20247.5
$ws_ARM.Range("L77").Value = 42847.855
$ws_ARM.Range("M77").Value = -15879.5
$ws_ARM.Range("N77").Value = -51583.855

# ARM row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
$ws_ARM.Range("H102").Value = 2026.4615
$ws_ARM.Range("I102").Value = 1819.5
$ws_ARM.Range("J102").Value = 2716.3333
$ws_ARM.Range("K102").Value = 1819.5
$ws_ARM.Range("L102").Value = 2716.3333
$ws_ARM.Range("M102").Value = -197.5
$ws_ARM.Range("N102").Value = -5960.3333

# ARM row 122: Haste for High Durium | High Durium Nugget
$ws_ARM.Range("H122").Value = 2092.3928
$ws_ARM.Range("I122").Value = 1944.8334
$ws_ARM.Range("K122").Value = 5834.5002
$ws_ARM.Range("M122").Value = -3384.5002

# ARM row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws_ARM.Range("H132").Value = 13205.096
$ws_ARM.Range("I132").Value = 10850.5
$ws_ARM.Range("K132").Value = 32551.5
$ws_ARM.Range("M132").Value = -30021.5

$ws_BSM = $wb.Worksheets.Item("BSM")
# BSM row 16: Port of Call: Ul'dah | Bronze Knuckles
$ws_BSM.Range("H16").Value = 814
$ws_BSM.Range("I16").Value = 814
$ws_BSM.Range("J16").Value = 0
$ws_BSM.Range("K16").Value = 814
$ws_BSM.Range("L16").Value = 0
$ws_BSM.Range("M16").ClearContents()
$ws_BSM.Range("N16").Value = -644

# BSM row 86: Through Thick and Thin | Adamantite Nugget
$ws_BSM.Range("H86").Value = 479180.2
$ws_BSM.Range("I86").Value = 716713.2
$ws_BSM.Range("K86").Value = 716713.2
$ws_BSM.Range("M86").Value = -715590.2

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws_BSM.Range("H89").Value = 479180.2
$ws_BSM.Range("I89").Value = 716713.2
$ws_BSM.Range("K89").Value = 3583566
$ws_BSM.Range("M89").Value = -3577950

# BSM row 94: High Steal | High Steel Nugget
$ws_BSM.Range("H94").Value = 550.4
$ws_BSM.Range("I94").Value = 167.48148
$ws_BSM.Range("K94").Value = 167.48148
$ws_BSM.Range("M94").Value = 283.51852

# BSM row 105: Ingot to Wing It | Molybdenum Ingot
$ws_BSM.Range("H105").Value = 2904.9033
$ws_BSM.Range("I105").Value = 3062.6843
$ws_BSM.Range("K105").Value = 3062.6843
$ws_BSM.Range("M105").Value = -1315.6843

# BSM row 107: The Gold Experience | Deepgold Nugget
$ws_BSM.Range("H107").Value = 1630.5
$ws_BSM.Range("I107").Value = 1651.2222
$ws_BSM.Range("K107").Value = 1651.2222
$ws_BSM.Range("M107").Value = 268.7778000000001

$ws_CRP = $wb.Worksheets.Item("CRP")
# CRP row 23: Nothing to Hide | Ash Mask (Lapis Lazuli)
$ws_CRP.Range("H23").Value = 32990
$ws_CRP.Range("J23").Value = 32990
$ws_CRP.Range("L23").Value = 32990
$ws_CRP.Range("N23").Value = -33470

# CRP row 27: Behind the Mask | Ash Mask (Lapis Lazuli)
$ws_CRP.Range("H27").Value = 32990
$ws_CRP.Range("J27").Value = 32990
$ws_CRP.Range("L27").Value = 32990
$ws_CRP.Range("N27").Value = -33374

# CRP row 31: Wall Not Found | Walnut Lumber
$ws_CRP.Range("H31").Value = 2465.4688
$ws_CRP.Range("J31").Value = 4078
$ws_CRP.Range("L31").Value = 4078
$ws_CRP.Range("N31").Value = -4668

# CRP row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws_CRP.Range("H34").Value = 2465.4688
$ws_CRP.Range("J34").Value = 4078
$ws_CRP.Range("L34").Value = 4078
$ws_CRP.Range("N34").Value = -4482

# CRP row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws_CRP.Range("H58").Value = 3717.0476
$ws_CRP.Range("I58").Value = 1275.3
$ws_CRP.Range("J58").Value = 5936.8184
$ws_CRP.Range("K58").Value = 1275.3
$ws_CRP.Range("L58").Value = 5936.8184
$ws_CRP.Range("M58").Value = -1072.3
$ws_CRP.Range("N58").Value = -6342.8184

# CRP row 99: O Pine | Pine Lumber
$ws_CRP.Range("H99").Value = 9762.486999999999
$ws_CRP.Range("I99").Value = 6783.1763
$ws_CRP.Range("J99").Value = 12064.682
$ws_CRP.Range("K99").Value = 6783.1763
$ws_CRP.Range("L99").Value = 12064.682
$ws_CRP.Range("M99").Value = -5285.1763
$ws_CRP.Range("N99").Value = -15060.682

# CRP row 107: Built to Last | White Oak Lumber
$ws_CRP.Range("H107").Value = 370.6389
$ws_CRP.Range("I107").Value = 286.41666
$ws_CRP.Range("J107").Value = 539.0833
$ws_CRP.Range("K107").Value = 286.41666
$ws_CRP.Range("L107").Value = 539.0833
$ws_CRP.Range("M107").Value = 1633.58334
$ws_CRP.Range("N107").Value = -4379.0833

# CRP row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws_CRP.Range("H122").Value = 2818.2222
$ws_CRP.Range("I122").Value = 2895.2
$ws_CRP.Range("K122").Value = 8685.599999999999
$ws_CRP.Range("M122").Value = -6235.599999999999

# CRP row 126: A Better Conductor | Red Pine Lumber
$ws_CRP.Range("H126").Value = 9762.486999999999
$ws_CRP.Range("I126").Value = 6783.1763
$ws_CRP.Range("J126").Value = 12064.682
$ws_CRP.Range("K126").Value = 20349.5289
$ws_CRP.Range("L126").Value = 36194.046
$ws_CRP.Range("M126").Value = -17879.5289
$ws_CRP.Range("N126").Value = -41134.046

# CRP row 132: Hull Lotta Damage | Ginseng Lumber
$ws_CRP.Range("H132").Value = 35628.74
$ws_CRP.Range("J132").Value = 42515
$ws_CRP.Range("L132").Value = 127545
$ws_CRP.Range("N132").Value = -132605

# CRP row 136: Turali Quality | Dark Mahogany Lumber
$ws_CRP.Range("H136").Value = 3717.0476
$ws_CRP.Range("I136").Value = 1275.3
$ws_CRP.Range("J136").Value = 5936.8184
$ws_CRP.Range("K136").Value = 3825.9
$ws_CRP.Range("L136").Value = 17810.4552
$ws_CRP.Range("M136").Value = -1275.9
$ws_CRP.Range("N136").Value = -22910.4552

$ws_GSM = $wb.Worksheets.Item("GSM")
# GSM row 104: Speak Softly and Carry a Metal Rod | Palladium Rod
$ws_GSM.Range("H104").Value = 85000
$ws_GSM.Range("J104").Value = 85000
$ws_GSM.Range("L104").Value = 85000
$ws_GSM.Range("N104").Value = -91988

# GSM row 126: Gold Rush Order | Phrygian Gold Ingot
$ws_GSM.Range("H126").Value = 2448.6
$ws_GSM.Range("I126").Value = 2439.2
$ws_GSM.Range("J126").Value = 2458
$ws_GSM.Range("K126").Value = 7317.599999999999
$ws_GSM.Range("L126").Value = 7374
$ws_GSM.Range("M126").Value = -4847.599999999999
$ws_GSM.Range("N126").Value = -12314

# GSM row 139: Ringing Gratitude | White Gold Ring of Healing
$ws_GSM.Range("H139").Value = 59326
$ws_GSM.Range("J139").Value = 59326
$ws_GSM.Range("L139").Value = 59326
$ws_GSM.Range("N139").Value = -69606

# GSM row 141: Mask Maker | Black Star Mask of Casting
$ws_GSM.Range("H141").Value = 100000
$ws_GSM.Range("J141").Value = 100000
$ws_GSM.Range("L141").Value = 100000
$ws_GSM.Range("N141").Value = -110360

$ws_LTW = $wb.Worksheets.Item("LTW")
# LTW row 22: Skin off Their Backs | Aldgoat Leather
$ws_LTW.Range("H22").Value = 2961.7585
$ws_LTW.Range("J22").Value = 4865.375
$ws_LTW.Range("L22").Value = 4865.375
$ws_LTW.Range("N22").Value = -5455.375

# LTW row 27: Fire and Hide | Aldgoat Leather
$ws_LTW.Range("H27").Value = 2961.7585
$ws_LTW.Range("J27").Value = 4865.375
$ws_LTW.Range("L27").Value = 4865.375
$ws_LTW.Range("N27").Value = -5079.375

# LTW row 46: Supply Side Logic | Boar Leather
$ws_LTW.Range("H46").Value = 1409.28
$ws_LTW.Range("I46").Value = 1066.3334
$ws_LTW.Range("K46").Value = 1066.3334
$ws_LTW.Range("M46").Value = -878.3334

# LTW row 55: It's Not a Job, It's a Calling | Peiste Leather
$ws_LTW.Range("H55").Value = 1055.4667
$ws_LTW.Range("I55").Value = 1163.5
$ws_LTW.Range("J55").Value = 932
$ws_LTW.Range("K55").Value = 1163.5
$ws_LTW.Range("L55").Value = 932
$ws_LTW.Range("M55").Value = -990.5
$ws_LTW.Range("N55").Value = -1278

# LTW row 100: Tiger in the Sack | Tiger Leather
$ws_LTW.Range("H100").Value = 5333
$ws_LTW.Range("I100").Value = 3999.5
$ws_LTW.Range("K100").Value = 3999.5
$ws_LTW.Range("M100").Value = -3458.5

# LTW row 132: Tenets of Tanning | Silver Lobo Leather
$ws_LTW.Range("H132").Value = 5360.125
$ws_LTW.Range("I132").Value = 5566.846
$ws_LTW.Range("K132").Value = 16700.538
$ws_LTW.Range("M132").Value = -14170.538

$ws_WVR = $wb.Worksheets.Item("WVR")
# WVR row 21: Don't Trew So Hard | Initiate's Slops
$ws_WVR.Range("H21").Value = 24006.4
$ws_WVR.Range("I21").Value = 19999
$ws_WVR.Range("J21").Value = 25008.25
$ws_WVR.Range("K21").Value = 19999
$ws_WVR.Range("L21").Value = 25008.25
$ws_WVR.Range("M21").Value = -19764
$ws_WVR.Range("N21").Value = -25478.25

# WVR row 35: Pantser Corps | Initiate's Slops
$ws_WVR.Range("H35").Value = 24006.4
$ws_WVR.Range("I35").Value = 19999
$ws_WVR.Range("J35").Value = 25008.25
$ws_WVR.Range("K35").Value = 19999
$ws_WVR.Range("L35").Value = 25008.25
$ws_WVR.Range("M35").Value = -19709
$ws_WVR.Range("N35").Value = -25588.25

# WVR row 42: Put on Your Party Pants | Velveteen Gaskins
$ws_WVR.Range("H42").Value = 99999
$ws_WVR.Range("I42").Value = 0
$ws_WVR.Range("K42").Value = 0
$ws_WVR.Range("M42").ClearContents()

# WVR row 107: Flax Wax | Bright Linen Yarn
$ws_WVR.Range("H107").Value = 948.7143
$ws_WVR.Range("I107").Value = 642.35297
$ws_WVR.Range("K107").Value = 1927.05891
$ws_WVR.Range("M107").Value = -7.058910000000196

# WVR row 126: A Polished Purchase | Snow Linen
$ws_WVR.Range("H126").Value = 6563.9614
$ws_WVR.Range("I126").Value = 6775.4287
$ws_WVR.Range("K126").Value = 20326.2861
$ws_WVR.Range("M126").Value = -17856.2861

# WVR row 132: Comfy Cabins | Snow Cotton Cloth
$ws_WVR.Range("H132").Value = 18518.5
$ws_WVR.Range("I132").Value = 11312.862
$ws_WVR.Range("K132").Value = 33938.586
$ws_WVR.Range("M132").Value = -31408.586
